# Auto-generated edit script: updates Leve profit calculation values
# across multiple sheets, reflecting refreshed market-board price data.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 70
$ws.Range("H70").Value = 1233.9642
$ws.Range("I70").Value = 987.2414
$ws.Range("J70").Value = 1498.963
$ws.Range("K70").Value = 2961.7242
$ws.Range("L70").Value = 4496.889
$ws.Range("M70").Value = -2691.7242
$ws.Range("N70").Value = -5036.889
# Row 73
$ws.Range("H73").Value = 1233.9642
$ws.Range("I73").Value = 987.2414
$ws.Range("J73").Value = 1498.963
$ws.Range("K73").Value = 2961.7242
$ws.Range("L73").Value = 4496.889
$ws.Range("M73").Value = -2025.7242
$ws.Range("N73").Value = -6368.889
# Row 112
$ws.Range("H112").Value = 3323.2942
$ws.Range("J112").Value = 3646.4
$ws.Range("L112").Value = 10939.2
$ws.Range("N112").Value = -13155.2
# Row 125
$ws.Range("H125").Value = 5228.2104
$ws.Range("I125").Value = 5000
$ws.Range("J125").Value = 5255.0586
$ws.Range("K125").Value = 45000
$ws.Range("L125").Value = 47295.52740000001
$ws.Range("M125").Value = -42540
$ws.Range("N125").Value = -52215.52740000001
# Row 129
$ws.Range("H129").Value = 1504.8077
$ws.Range("I129").Value = 777.5
$ws.Range("J129").Value = 1565.4166
$ws.Range("K129").Value = 2332.5
$ws.Range("L129").Value = 4696.2498
$ws.Range("M129").Value = 2667.5
$ws.Range("N129").Value = -14696.2498
# Row 132
$ws.Range("H132").Value = 6424.3184
$ws.Range("I132").Value = 4971.8184
$ws.Range("J132").Value = 10781.818
$ws.Range("K132").Value = 14915.4552
$ws.Range("L132").Value = 32345.454
$ws.Range("M132").Value = -12385.4552
$ws.Range("N132").Value = -37405.454
# Row 138
$ws.Range("H138").Value = 2190.4807
$ws.Range("I138").Value = 1656.2
$ws.Range("J138").Value = 2685.1853
$ws.Range("K138").Value = 4968.6
$ws.Range("L138").Value = 8055.5559
$ws.Range("M138").Value = 171.3999999999996
$ws.Range("N138").Value = -18335.5559

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 3839702.5
$ws.Range("I32").Value = 3955.7534
$ws.Range("K32").Value = 3955.7534
$ws.Range("M32").Value = -3668.7534
# Row 74
$ws.Range("H74").Value = 14707199
$ws.Range("I74").Value = 20834722
$ws.Range("J74").Value = 1145.2
$ws.Range("K74").Value = 20834722
$ws.Range("L74").Value = 1145.2
$ws.Range("M74").Value = -20833848
$ws.Range("N74").Value = -2893.2
# Row 77
$ws.Range("H77").Value = 14707199
$ws.Range("I77").Value = 20834722
$ws.Range("J77").Value = 1145.2
$ws.Range("K77").Value = 104173610
$ws.Range("L77").Value = 5726
$ws.Range("M77").Value = -104169242
$ws.Range("N77").Value = -14462
# Row 122
$ws.Range("H122").Value = 3919.7
$ws.Range("I122").Value = 1576
$ws.Range("K122").Value = 4728
$ws.Range("M122").Value = -2278

$ws = $wb.Worksheets.Item("BSM")
# Row 94
$ws.Range("H94").Value = 455.48
$ws.Range("I94").Value = 382.13635
$ws.Range("J94").Value = 993.3333
$ws.Range("K94").Value = 382.13635
$ws.Range("L94").Value = 993.3333
$ws.Range("M94").Value = 68.86365000000001
$ws.Range("N94").Value = -1895.3333
# Row 134
$ws.Range("H134").Value = 5807.081
$ws.Range("I134").Value = 2176.0588
$ws.Range("J134").Value = 8893.450000000001
$ws.Range("K134").Value = 6528.176399999999
$ws.Range("L134").Value = 26680.35
$ws.Range("M134").Value = -3993.176399999999
$ws.Range("N134").Value = -31750.35

$ws = $wb.Worksheets.Item("CRP")
# Row 58
$ws.Range("H58").Value = 2008273.5
$ws.Range("I58").Value = 1602.6
$ws.Range("J58").Value = 5018280
$ws.Range("K58").Value = 1602.6
$ws.Range("L58").Value = 5018280
$ws.Range("M58").Value = -1399.6
$ws.Range("N58").Value = -5018686
# Row 94
$ws.Range("H94").Value = 1886
$ws.Range("I94").Value = 1860.5
$ws.Range("K94").Value = 1860.5
$ws.Range("M94").Value = -1409.5
# Row 107
$ws.Range("H107").Value = 2092.375
$ws.Range("I107").Value = 943.1111
$ws.Range("J107").Value = 3570
$ws.Range("K107").Value = 943.1111
$ws.Range("L107").Value = 3570
$ws.Range("M107").Value = 976.8889
$ws.Range("N107").Value = -7410
# Row 134
$ws.Range("H134").Value = 3200.12
$ws.Range("I134").Value = 1808.6364
$ws.Range("J134").Value = 4293.4287
$ws.Range("K134").Value = 5425.9092
$ws.Range("L134").Value = 12880.2861
$ws.Range("M134").Value = -2890.9092
$ws.Range("N134").Value = -17950.2861
# Row 136
$ws.Range("H136").Value = 2008273.5
$ws.Range("I136").Value = 1602.6
$ws.Range("J136").Value = 5018280
$ws.Range("K136").Value = 4807.799999999999
$ws.Range("L136").Value = 15054840
$ws.Range("M136").Value = -2257.799999999999
$ws.Range("N136").Value = -15059940

$ws = $wb.Worksheets.Item("CUL")
# Row 68
$ws.Range("H68").Value = 948026.6
$ws.Range("I68").Value = 976.2646999999999
$ws.Range("J68").Value = 3248006
$ws.Range("K68").Value = 2928.7941
$ws.Range("L68").Value = 9744018
$ws.Range("M68").Value = -2117.7941
$ws.Range("N68").Value = -9745640
# Row 71
$ws.Range("H71").Value = 948026.6
$ws.Range("I71").Value = 976.2646999999999
$ws.Range("J71").Value = 3248006
$ws.Range("K71").Value = 8786.382299999999
$ws.Range("L71").Value = 29232054
$ws.Range("M71").Value = -4730.382299999999
$ws.Range("N71").Value = -29240166
# Row 105
$ws.Range("H105").Value = 6081.8125
$ws.Range("J105").Value = 6287.2666
$ws.Range("L105").Value = 18861.7998
$ws.Range("N105").Value = -24103.7998
# Row 131
$ws.Range("H131").Value = 291607.25
$ws.Range("I131").Value = 521.8125
$ws.Range("J131").Value = 464102.34
$ws.Range("K131").Value = 1565.4375
$ws.Range("L131").Value = 1392307.02
$ws.Range("M131").Value = 3474.5625
$ws.Range("N131").Value = -1402387.02

$ws = $wb.Worksheets.Item("GSM")
# Row 132
$ws.Range("H132").Value = 1553.5555
$ws.Range("I132").Value = 1246.625
$ws.Range("J132").Value = 2000
$ws.Range("K132").Value = 3739.875
$ws.Range("L132").Value = 6000
$ws.Range("M132").Value = -1209.875
$ws.Range("N132").Value = -11060

$ws = $wb.Worksheets.Item("LTW")
# Row 29
$ws.Range("H29").Value = 38259
$ws.Range("I29").Value = 5000
$ws.Range("J29").Value = 49345.332
$ws.Range("K29").Value = 5000
$ws.Range("L29").Value = 49345.332
$ws.Range("M29").Value = -4705
$ws.Range("N29").Value = -49935.332
# Row 33
$ws.Range("H33").Value = 30000
$ws.Range("I33").Value = 30000
$ws.Range("J33").Value = 0
$ws.Range("K33").Value = 30000
$ws.Range("L33").Value = 0
$ws.Range("M33").Value = -29710
$ws.Range("N33").ClearContents()
# Row 122
$ws.Range("H122").Value = 9546.522999999999
$ws.Range("I122").Value = 11552.833
$ws.Range("J122").Value = 6871.4443
$ws.Range("K122").Value = 34658.499
$ws.Range("L122").Value = 20614.3329
$ws.Range("M122").Value = -32208.499
$ws.Range("N122").Value = -25514.3329

$ws = $wb.Worksheets.Item("WVR")
# Row 32
$ws.Range("H32").Value = 2995
$ws.Range("I32").Value = 2995
$ws.Range("K32").Value = 2995
$ws.Range("M32").Value = -2678
# Row 132
$ws.Range("H132").Value = 2056.125
$ws.Range("I132").Value = 1355.4814
$ws.Range("J132").Value = 2956.9524
$ws.Range("K132").Value = 4066.4442
$ws.Range("L132").Value = 8870.8572
$ws.Range("M132").Value = -1536.4442
$ws.Range("N132").Value = -13930.8572
